# ------------------------------------------------------------------------
# Reproduces the "compare_old.xlsx" edit:
#  - Sheet1 gets a new "VQSR" column of PASS/FAIL values (col E) plus a
#    numeric p-value-ish column (col F) for a couple of rows.
#  - The active selection moves from F8 to A6.
#  - Page setup is switched to Letter/portrait-ish printing (paperSize 9).
#  - Two (unused) built-in "Hyperlink" / "Followed Hyperlink" cell styles
#    end up registered in the workbook style table (a side effect of some
#    hyperlinked data having been pasted in and then overwritten with
#    plain values), without leaving any real hyperlink or cell-level
#    style behind on the sheet itself.
# ------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New / changed cell values -------------------------------------------------
$ws.Range("E2").Value = "VQSR"
$ws.Range("F2").Value = 0.3

$ws.Range("E3").Value = "PASS"

$ws.Range("E4").Value = "FAIL"

$ws.Range("E5").Value = "PASS"
$ws.Range("F5").Value = 0.004

$ws.Range("E6").Value = "PASS"
$ws.Range("F6").Value = 0.0003

# --- Register the "Hyperlink" / "Followed Hyperlink" built-in cell styles ------
# (done on a throw-away worksheet so the visible sheet / its used range is
# left completely untouched, matching the fact that no cell on Sheet1 ends
# up referencing these styles)
$tmpWs = $wb.Worksheets.Add()
$tmpCell1 = $tmpWs.Range("A1")
$tmpCell1.Hyperlinks.Add($tmpCell1, "http://example.com")
$tmpCell2 = $tmpWs.Range("A2")
$tmpCell2.Style = "Followed Hyperlink"
$tmpWs.Delete()

# --- Page setup ------------------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# --- Window geometry --------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = 720
$win.Top = 220
$win.Width = 14300
$win.Height = 15240

# --- Active selection --------------------------------------------------------
$ws.Range("A6").Select()
